$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (ECs target) with new TPM-derived values
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.08037899999999999
$ws.Range("N2").Value = 0.241137
$ws.Range("O2").Value = 0.01215705881223039
$ws.Range("P2").Value = 0.01215705881223039
$ws.Range("Q2").Value = 0.130360966398
$ws.Range("R2").Value = 1.173248697582
$ws.Range("S2").Value = 0.01215705881223039
$ws.Range("T2").Value = 0.01215705881223039

# Update row 3 (FAPs target) - only O,P,S,T change
$ws.Range("O3").Value = 0.2935162100923598
$ws.Range("P3").Value = 0.2935162100923598
$ws.Range("S3").Value = 0.2935162100923598
$ws.Range("T3").Value = 0.2935162100923598

# Update row 4 (MuSCs target) with new TPM-derived values
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.590689999999999
$ws.Range("N4").Value = 13.77207
$ws.Range("O4").Value = 0.6943267310954097
$ws.Range("P4").Value = 0.6943267310954098
$ws.Range("Q4").Value = 7.445312641779998
$ws.Range("R4").Value = 67.00781377601999
$ws.Range("S4").Value = 0.6943267310954097
$ws.Range("T4").Value = 0.6943267310954098

# Remove row 5 entirely (Resolving-Mac target cluster row no longer present)
$ws.Rows.Item(5).Delete()
